$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B22").Value = 6294
$ws.Range("D22").Value = 5849586
$ws.Range("E22").Value = 929.3908484270734
$ws.Range("F22").Value = 8.349113444654854
$ws.Range("H22").Value = 27.21026626884366
